$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.682.61'
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").Value = '1.778.33'
$ws.Range("E3").Value = '  -0.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.76'
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.42%  '

$ws.Range("E7").Value = '  +3.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3439'
$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.73'
$ws.Range("E9").Value = '  -1.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.159'
$ws.Range("E10").Value = '  -3.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07543'
$ws.Range("E11").Value = '  +0.85%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.96'
$ws.Range("E12").Value = '  +5.79%  '

$ws.Range("B13").Value = 'BinanceUSD'
$ws.Range("C13").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.006'
$ws.Range("E13").Value = '  +0.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.457'
$ws.Range("E14").Value = '  -0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.185'
$ws.Range("E15").Value = '  +1.47%  '

$ws.Range("D16").Value = '1.781.12'
$ws.Range("E16").Value = '  -0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001089'
$ws.Range("E17").Value = '  -0.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06712'
$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.66'
$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.72'
$ws.Range("E21").Value = '  +2.42%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.535'
$ws.Range("E22").Value = '  -1.62%  '

$ws.Range("D23").Value = '27.701.28'
$ws.Range("E23").Value = '  +1.79%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.34'
$ws.Range("E24").Value = '  -0.45%  '

$ws.Range("E25").Value = '  -1.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.491'
$ws.Range("E26").Value = '  -1.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.511'
$ws.Range("E27").Value = '  -1.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.09'
$ws.Range("E28").Value = '  -1.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '155.03'
$ws.Range("E29").Value = '  +0.96%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '137.08'
$ws.Range("E30").Value = '  +2.27%  '

$ws.Range("D31").Value = '1.983.14'
$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.284'
$ws.Range("E32").Value = '  +3.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.980'
$ws.Range("E33").Value = '  -0.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08936'
$ws.Range("E34").Value = '  +3.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.02'
$ws.Range("E35").Value = '  -1.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02472'
$ws.Range("E36").Value = '  +5.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.488'
$ws.Range("E37").Value = '  +0.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6918'
$ws.Range("E38").Value = '  -0.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06464'
$ws.Range("E39").Value = '  +1.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2229'
$ws.Range("E40").Value = '  +1.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.581'
$ws.Range("E41").Value = '  -4.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.262'
$ws.Range("E42").Value = '  +1.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.541'
$ws.Range("E43").Value = '  -2.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.42'
$ws.Range("E44").Value = '  +0.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.004'
$ws.Range("E45").Value = '  +0.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6389'
$ws.Range("E46").Value = '  -1.56%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.863'
$ws.Range("E47").Value = '  +0.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.153'
$ws.Range("E48").Value = '  +0.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '133.50'
$ws.Range("E49").Value = '  +3.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07487'
$ws.Range("E50").Value = '  +5.06%  '

$ws.Range("B51").Value = 'Tezos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.259'
$ws.Range("E51").Value = '  +4.23%  '
